$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Row 17 ("repaymentstrategy") value changes from "RBI (India)" to the new
# overdue/due fee scenario string as part of adding periodic & upfront
# related scenarios.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Reflect the new active cell / selection in the sheet view.
$ws.Range("B17").Select()
